$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(13)
$ws.Range("L3").Value = "WFH"
